$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries (shared-string table reorder reflected as new country per row) ---
# Japon, Polonia, Emiratos Arabes Unidos, Rumania -> Japon, Emiratos Arabes Unidos, Polonia, Rumania
$ws.Range("A35").Value = 'Emiratos Arabes Unidos'
$ws.Range("A36").Value = 'Polonia'
$ws.Range("A111").Value = 'Mali'
$ws.Range("A112").Value = 'Mayotte'
$ws.Range("A113").Value = 'Sudan'
$ws.Range("A114").Value = 'Maldivas'
$ws.Range("A138").Value = 'Cabo Verde'
$ws.Range("A139").Value = 'Liberia'
$ws.Range("A140").Value = 'Guadalupe'
$ws.Range("A141").Value = 'Birmania'
$ws.Range("A142").Value = 'Gibraltar'
$ws.Range("A143").Value = 'Brunei'
$ws.Range("A144").Value = 'Madagascar'
$ws.Range("A145").Value = 'Etiopia'
$ws.Range("A146").Value = 'Guayana Francesa'
$ws.Range("A147").Value = 'Togo'
$ws.Range("A193").Value = 'San Vicente y las Granadinas'
$ws.Range("A194").Value = 'Namibia'
$ws.Range("A217").Value = 'San Pedro y Miquelon'
$ws.Range("A218").Value = 'Comoras'

# --- Updated case numbers (Casos totales/Nuevos/Activos/Recuperados/Criticos/MuertesHoy/Muertes) ---
# Row 4
$ws.Range("B4").Value = 1146766
$ws.Range("C4").Value = 15736
$ws.Range("E4").Value = 918038
$ws.Range("G4").Value = 868
$ws.Range("H4").Value = 66621
# Row 9
$ws.Range("F9").Value = 2105
# Row 15
$ws.Range("B15").Value = 56580
$ws.Range("C15").Value = 1519
$ws.Range("E15").Value = 29704
$ws.Range("G15").Value = 169
$ws.Range("H15").Value = 3560
# Row 35
$ws.Range("B35").Value = 13599
$ws.Range("C35").Value = 561
$ws.Range("D35").Value = 2664
$ws.Range("E35").Value = 10816
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 8
$ws.Range("H35").Value = 119
# Row 36
$ws.Range("B36").Value = 13375
$ws.Range("C36").Value = 270
$ws.Range("D36").Value = 3762
$ws.Range("E36").Value = 8949
$ws.Range("F36").Value = 160
$ws.Range("G36").Value = 13
$ws.Range("H36").Value = 664
# Row 45
$ws.Range("B45").Value = 7801
$ws.Range("C45").Value = 18
$ws.Range("E45").Value = 7559
# Row 111
$ws.Range("B111").Value = 544
$ws.Range("C111").Value = 36
$ws.Range("D111").Value = 206
$ws.Range("E111").Value = 312
$ws.Range("F111").Value = 0
$ws.Range("H111").Value = 26
# Row 112
$ws.Range("B112").Value = 539
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 235
$ws.Range("E112").Value = 300
$ws.Range("F112").Value = 4
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 4
# Row 113
$ws.Range("B113").Value = 533
$ws.Range("C113").Value = 91
$ws.Range("D113").Value = 46
$ws.Range("E113").Value = 451
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 5
$ws.Range("H113").Value = 36
# Row 114
$ws.Range("B114").Value = 514
$ws.Range("C114").Value = 23
$ws.Range("D114").Value = 17
$ws.Range("E114").Value = 496
$ws.Range("F114").Value = 2
$ws.Range("H114").Value = 1
# Row 138
$ws.Range("C138").Value = 30
$ws.Range("D138").Value = 18
$ws.Range("E138").Value = 132
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 2
# Row 139
$ws.Range("D139").Value = 45
$ws.Range("E139").Value = 89
$ws.Range("F139").Value = 0
$ws.Range("H139").Value = 18
# Row 140
$ws.Range("B140").Value = 152
$ws.Range("D140").Value = 95
$ws.Range("E140").Value = 45
$ws.Range("F140").Value = 6
$ws.Range("H140").Value = 12
# Row 141
$ws.Range("B141").Value = 151
$ws.Range("D141").Value = 37
$ws.Range("E141").Value = 108
$ws.Range("H141").Value = 6
# Row 142
$ws.Range("B142").Value = 144
$ws.Range("D142").Value = 131
$ws.Range("E142").Value = 13
$ws.Range("F142").Value = 0
$ws.Range("H142").Value = 0
# Row 143
$ws.Range("B143").Value = 138
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 126
$ws.Range("E143").Value = 11
$ws.Range("F143").Value = 2
$ws.Range("H143").Value = 1
# Row 144
$ws.Range("B144").Value = 135
$ws.Range("C144").Value = 3
$ws.Range("D144").Value = 97
$ws.Range("E144").Value = 38
$ws.Range("F144").Value = 1
$ws.Range("H144").Value = 0
# Row 145
$ws.Range("B145").Value = 133
$ws.Range("D145").Value = 69
$ws.Range("E145").Value = 61
$ws.Range("F145").Value = 0
$ws.Range("H145").Value = 3
# Row 146
$ws.Range("B146").Value = 128
$ws.Range("D146").Value = 98
$ws.Range("E146").Value = 29
$ws.Range("F146").Value = 2
$ws.Range("H146").Value = 1
# Row 147
$ws.Range("B147").Value = 123
$ws.Range("D147").Value = 66
$ws.Range("E147").Value = 48
$ws.Range("H147").Value = 9
